$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 29; $row++) {
    $ws.Range("BF$row").Value = "2020-12-24"
}
